$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new test row, following the same pattern as row 13 but referencing A14/B14 etc.
$ws.Range("A14").Formula = "=2.5*B4"
$ws.Range("B14").Value = 1

$ws.Range("C14").Formula = "=0.04*A14 *B14"
$ws.Range("D14").Formula = "=0.045*A14*B14"
$ws.Range("E14").Formula = "=0.04*A14*B14"
$ws.Range("F14").Formula = "=0.08*A14*B14"
$ws.Range("G14").Formula = "=A14-C14/B14-E14/B14-H14/B14+IF(A14/`$B`$2 <= 2, `$A`$2, 0)"
$ws.Range("H14").Formula = "=A14*L14*B14"
$ws.Range("I14").Formula = "=A14-C14/B14-E14/B14-H14/B14"
$ws.Range("J14").Formula = "=(I14-I14*0.25)*B14"
$ws.Range("K14").Formula = "=B14 * IF(AND(M14 > 0, M14 <= 95),0,IF(AND(M14 > 95, M14 <= 150),((M14 - 95)*`$C`$2)*0.19,IF(AND(M14 > 150, M14 <= 360),((M14 - 150)*`$C`$2)*0.28 + 10*`$C`$2,((M14 - 360)*`$C`$2)*0.33 + 69*`$C`$2)))"
$ws.Range("L14").Formula = "=IF(A14/`$B`$2<4,0,IF(AND(A14/`$B`$2>=4,A14/`$B`$2<16),`$D`$2,IF(AND(A14/`$B`$2>=16,A14/`$B`$2<17),`$E`$2,IF(AND(A14/`$B`$2>=17,A14/`$B`$2<18),`$F`$2,IF(AND(A14/`$B`$2>=18,A14/`$B`$2<19),`$G`$2,IF(AND(A14/`$B`$2>=19,A14/`$B`$2<20),`$H`$2,`$I`$2))))))"
$ws.Range("M14").Formula = "=(J14/B14)/`$C`$2"
$ws.Range("N14").Formula = "=G14 * B14-(K14/B14)"

# Apply the number format to I14 (matches style index 1 used by I6:I13) after
# all sibling cells in the row have been created, to avoid the format
# bleeding into cells created afterwards.
$ws.Range("I14").NumberFormat = "#,##0"

$ws.Range("K11").Select()
